# Commit: "renaming industry by site"
# Cell A7 on the only worksheet reads:
#   "Increase in temperature in the receiving water body due to industry
#    discharge" (+ a trailing bold-superscript "1" footnote marker)
# and must become:
#   "Increase in temperature in the receiving water body due to site
#    discharge" (+ the same trailing bold-superscript "1" footnote marker)
# No other cell's displayed text changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a7Text = "Increase in temperature in the receiving water body due to site discharge1"

$cellA7 = $ws.Range("A7")
$cellA7.Value = $a7Text

# Re-apply the trailing "1" footnote marker as bold superscript, matching
# the formatting already used for this reference mark.
$fontA7 = $cellA7.Characters($a7Text.Length, 1).Font()
$fontA7.Bold = $true
$fontA7.Superscript = $true
